# Update the Lunch Sponsor amount from 600 to 1000
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D10").Value = 1000

# Update the view/selection state to match the saved workbook view
$excel.ActiveWindow.ScrollRow = 6
$ws.Range("D11").Select()
